# Refactor timetable generation to apply consistent cell alignment and
# borders for improved visual presentation: center + middle-align + wrap
# every bordered content cell (including the bold header row), and fix up
# the Saturday (column G) block whose rows/merges had drifted by one row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$xlCenter = -4108

# ---------------------------------------------------------------------
# 1) Column G (Saturday) content + merge fix-up.
#    Before:
#      G15:G18 merged  -> "Lunch / Dress Up, Warm Up"
#      G19     standalone -> "Concert call time"
#      G20:G25 merged  -> "Lina Summer Camp of Music Students & Friends Concert"
#      G26:G28 merged  -> "After concert refreshment / (Maritime Museum)"
#    After:
#      G15:G19 merged  -> "Lunch / Dress Up, Warm Up"   (extended by 1 row)
#      G20     standalone -> "Concert call time"         (shifted down 1 row)
#      G21:G26 merged  -> "Lina Summer Camp ..."         (shifted down 1 row)
#      G27:G28 merged  -> "After concert refreshment..." (shifted down 1 row)
# ---------------------------------------------------------------------

# Unmerge the affected ranges first so every cell becomes individually
# addressable (Excel only lets you write the top-left/anchor cell of a
# merged range).
$ws.Range("G15:G18").UnMerge()
$ws.Range("G20:G25").UnMerge()
$ws.Range("G26:G28").UnMerge()

# Move the content down by one row.
$ws.Range("G20").Value2 = "Concert call time"
$ws.Range("G21").Value2 = "Lina Summer Camp of Music Students & Friends Concert"
$ws.Range("G27").Value2 = "After concert refreshment `n(Maritime Museum)"

# Clear the cells whose content moved away.
$ws.Range("G19").Value2 = $null
$ws.Range("G26").Value2 = $null

# Re-merge with the new, shifted extents.
$ws.Range("G15:G19").Merge()
$ws.Range("G21:G26").Merge()
$ws.Range("G27:G28").Merge()

# Merging/unmerging re-derives per-cell borders from the merged block's
# outer edge (interior edges lose their border), same as Excel's own
# "smart" merge border behaviour. This workbook always uses a plain thin
# box around every cell, merged or not, so restore the bottom edge (which
# re-unions it back to the uniform thin-box style used everywhere else).
$ws.Range("G15:G19").Borders.Item(4).LineStyle = 1
$ws.Range("G21:G26").Borders.Item(4).LineStyle = 1
$ws.Range("G27:G28").Borders.Item(4).LineStyle = 1

# ---------------------------------------------------------------------
# 2) Consistent alignment: every border-ringed content cell (the ones that
#    previously only had vertical-center + wrap, plus the bold header
#    cells) now also gets horizontal centering - center/middle/wrap.
# ---------------------------------------------------------------------

$alignRefs = @(
    "B3","C3","D3","E3","F3","G4",
    "B7","C7","D7","E7","F7","G7",
    "B11","C11","E11",
    "B15","C15","D15","E15","F15","G15",
    "B19","C19","D19","E19","F19",
    "B23","C23","D23","E23","F23",
    "B27","C27","D27","E27","F27",
    "G20","G21","G27",
    "B1","C1","D1","E1","F1","G1","A2"
)

foreach ($ref in $alignRefs) {
    $r = $ws.Range($ref)
    $r.HorizontalAlignment = $xlCenter
    $r.VerticalAlignment = $xlCenter
    $r.WrapText = $true
}
